# Weekly update: insert a new daily price record for
# "Vega Monumental Concepción - Cebollín" before the existing row 53,
# shifting all subsequent records (old rows 53-126) down by one row
# (they become rows 54-127), and populate the new row 53 with the
# latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; Excel shifts rows 53..126 down to 54..127,
# carrying along the existing "Fecha" (column D) number format (style id 2).
$ws.Rows.Item(53).Insert()

# Columns that stay constant across every record in this sheet (Mercado ID,
# Mercado, Región, Codreg, Categoría ID, Categoría, Variedad, Clasificación)
# are copied from the neighboring row so the new entry is consistent.
$ws.Cells.Item(53, 1).Value2 = $ws.Cells.Item(54, 1).Value2   # Mercado ID
$ws.Cells.Item(53, 2).Value2 = $ws.Cells.Item(54, 2).Value2   # Mercado
$ws.Cells.Item(53, 3).Value2 = $ws.Cells.Item(54, 3).Value2   # Región
$ws.Cells.Item(53, 5).Value2 = $ws.Cells.Item(54, 5).Value2   # Codreg
$ws.Cells.Item(53, 6).Value2 = $ws.Cells.Item(54, 6).Value2   # Categoría ID
$ws.Cells.Item(53, 7).Value2 = $ws.Cells.Item(54, 7).Value2   # Categoría
$ws.Cells.Item(53, 8).Value2 = $ws.Cells.Item(54, 8).Value2   # Variedad
$ws.Cells.Item(53, 18).Value2 = $ws.Cells.Item(54, 18).Value2 # Clasificación

# New observation's own data.
$ws.Cells.Item(53, 4).Value2 = 45128                          # Fecha
$ws.Cells.Item(53, 9).Value2 = "Primera"                      # Calidad
$ws.Cells.Item(53, 10).Value2 = 60                            # Volumen
$ws.Cells.Item(53, 11).Value2 = 5000                          # Precio mínimo
$ws.Cells.Item(53, 12).Value2 = 5000                          # Precio máximo
$ws.Cells.Item(53, 13).Value2 = 5000                          # Precio promedio ponderado
$ws.Cells.Item(53, 14).Value2 = "$/paquete 36 unidades"       # Unidad de comercialización
$ws.Cells.Item(53, 15).Value2 = "Región Metropolitana"        # Origen
$ws.Cells.Item(53, 16).Value2 = 139                           # Precio $/Kg
$ws.Cells.Item(53, 17).Value2 = 36                            # Kg o Unidades
